$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing text storage so that
# numeric-looking strings (e.g. "314.82", "3.41") are not silently
# re-typed as numbers by Excel's usual auto-detection. The NumberFormat
# is restored to the sheet's default ("Normal" style) right afterwards so
# no visible formatting change is left behind.
function Set-TextValue($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "57.670.96"
$ws.Range("E2").Value = "  -0.74%  "
# Row 3
Set-TextValue $ws "D3" "2.456.64"
$ws.Range("E3").Value = "  +0.06%  "
# Row 4
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
Set-TextValue $ws "D5" "511.42"
# Row 6
Set-TextValue $ws "D6" "133.96"
$ws.Range("E6").Value = "  +2.59%  "
# Row 7
Set-TextValue $ws "D7" "0.997"
$ws.Range("E7").Value = "  -0.17%  "
# Row 8
$ws.Range("E8").Value = "  -1.62%  "
# Row 9
Set-TextValue $ws "D9" "2.458.44"
$ws.Range("E9").Value = "  -0.08%  "
# Row 10
Set-TextValue $ws "D10" "0.0979"
$ws.Range("E10").Value = "  -0.19%  "
# Row 11
$ws.Range("E11").Value = "  -0.83%  "
# Row 12
$ws.Range("E12").Value = "  -0.54%  "
# Row 13
$ws.Range("E13").Value = "  -7.38%  "
# Row 14
Set-TextValue $ws "D14" "2.889.67"
$ws.Range("E14").Value = "  -0.04%  "
# Row 15
Set-TextValue $ws "D15" "57.669.43"
$ws.Range("E15").Value = "  -0.54%  "
# Row 16
Set-TextValue $ws "D16" "21.96"
$ws.Range("E16").Value = "  +1.19%  "
# Row 17
$ws.Range("E17").Value = "  +0.82%  "
# Row 18
Set-TextValue $ws "D18" "2.439.34"
$ws.Range("E18").Value = "  -0.56%  "
# Row 19
Set-TextValue $ws "D19" "10.32"
$ws.Range("E19").Value = "  -1.20%  "
# Row 20
Set-TextValue $ws "D20" "4.13"
$ws.Range("E20").Value = "  -0.01%  "
# Row 21
Set-TextValue $ws "D21" "314.82"
# Row 22
$ws.Range("E22").Value = "  +3.80%  "
# Row 23
$ws.Range("E23").Value = "  +0.05%  "
# Row 24
$ws.Range("E24").Value = "  -2.17%  "
# Row 25
Set-TextValue $ws "D25" "65.26"
$ws.Range("E25").Value = "  +0.01%  "
# Row 26
Set-TextValue $ws "D26" "0.998"
$ws.Range("E26").Value = "  -0.15%  "
# Row 27
$ws.Range("E27").Value = "  -1.13%  "
# Row 28
$ws.Range("E28").Value = "  -6.10%  "
# Row 29
$ws.Range("E29").Value = "  +3.92%  "
# Row 30
Set-TextValue $ws "D30" "173.23"
$ws.Range("E30").Value = "  -0.36%  "
# Row 31
$ws.Range("E31").Value = "  -0.55%  "
# Row 32
Set-TextValue $ws "D32" "1.70"
$ws.Range("E32").Value = "  -0.07%  "
# Row 33
Set-TextValue $ws "D33" "6.18"
$ws.Range("E33").Value = "  +0.18%  "
# Row 34
$ws.Range("E34").Value = "  +0.33%  "
# Row 35
$ws.Range("E35").Value = "  +0.03%  "
# Row 36
Set-TextValue $ws "D36" "0.996"
$ws.Range("E36").Value = "  -0.12%  "
# Row 37
Set-TextValue $ws "D37" "18.06"
$ws.Range("E37").Value = "  +0.85%  "
# Row 38
$ws.Range("E38").Value = "  +5.17%  "
# Row 39
Set-TextValue $ws "D39" "3.87"
$ws.Range("E39").Value = "  +1.84%  "
# Row 40
$ws.Range("E40").Value = "  +1.05%  "
# Row 41
$ws.Range("E41").Value = "  +0.65%  "
# Row 42
Set-TextValue $ws "D42" "0.811"
$ws.Range("E42").Value = "  -0.44%  "
# Row 43
Set-TextValue $ws "D43" "136.51"
$ws.Range("E43").Value = "  +6.97%  "
# Row 44
Set-TextValue $ws "D44" "3.41"
$ws.Range("E44").Value = "  -0.21%  "
# Row 45
$ws.Range("E45").Value = "  +2.24%  "
# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D46" "0.577"
$ws.Range("E46").Value = "  -1.50%  "
# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D47" "256.89"
$ws.Range("E47").Value = "  -1.53%  "
# Row 48
Set-TextValue $ws "D48" "0.0918"
$ws.Range("E48").Value = "  -0.60%  "
# Row 49
Set-TextValue $ws "D49" "0.0494"
$ws.Range("E49").Value = "  +0.15%  "
# Row 50
Set-TextValue $ws "D50" "0.0215"
$ws.Range("E50").Value = "  +1.30%  "
# Row 51
Set-TextValue $ws "D51" "17.20"
$ws.Range("E51").Value = "  +0.88%  "
